# Auto-generated script applying scheduled market-data refresh to Sheets/Ifrit_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the rows
# whose underlying Universalis market data changed.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 721.9  # H29 (was 250)
$ws.Cells.Item(29, 9).Value = 800.75  # I29 (was 250)
$ws.Cells.Item(29, 10).Value = 669.3333  # J29 (was 0)
$ws.Cells.Item(29, 11).Value = 2402.25  # K29 (was 750)
$ws.Cells.Item(29, 12).Value = 2007.9999  # L29 (was 0)
$ws.Cells.Item(29, 13).Value = -2121.25  # M29 (was -469)
$ws.Cells.Item(29, 14).Value = -2569.9999  # N29 (was None)
$ws.Cells.Item(38, 8).Value = 67.07692  # H38 (was 117.92308)
$ws.Cells.Item(38, 9).Value = 67.07692  # I38 (was 69.416664)
$ws.Cells.Item(38, 10).Value = 0  # J38 (was 700)
$ws.Cells.Item(38, 11).Value = 201.23076  # K38 (was 208.249992)
$ws.Cells.Item(38, 12).Value = 0  # L38 (was 2100)
$ws.Cells.Item(38, 13).Value = 170.76924  # M38 (was 163.750008)
$ws.Cells.Item(38, 14).ClearContents() | Out-Null  # N38 (was -2844)
$ws.Cells.Item(40, 8).Value = 892.2308  # H40 (was 910)
$ws.Cells.Item(40, 9).Value = 892.2308  # I40 (was 910)
$ws.Cells.Item(40, 11).Value = 892.2308  # K40 (was 910)
$ws.Cells.Item(40, 13).Value = -717.2308  # M40 (was -735)
$ws.Cells.Item(43, 8).Value = 435  # H43 (was 500.25)
$ws.Cells.Item(43, 9).Value = 366.66666  # I43 (was 533.6667)
$ws.Cells.Item(43, 10).Value = 476  # J43 (was 400)
$ws.Cells.Item(43, 11).Value = 366.66666  # K43 (was 533.6667)
$ws.Cells.Item(43, 12).Value = 476  # L43 (was 400)
$ws.Cells.Item(43, 13).Value = -297.66666  # M43 (was -464.6667)
$ws.Cells.Item(43, 14).Value = -614  # N43 (was -538)
$ws.Cells.Item(58, 8).Value = 1391.3334  # H58 (was 1497.4546)
$ws.Cells.Item(58, 10).Value = 2000  # J58 (was 1993.75)
$ws.Cells.Item(58, 12).Value = 6000  # L58 (was 5981.25)
$ws.Cells.Item(58, 14).Value = -6300  # N58 (was -6281.25)
$ws.Cells.Item(86, 8).Value = 2200.7693  # H86 (was 2662.158)
$ws.Cells.Item(86, 9).Value = 1110.7333  # I86 (was 1335.125)
$ws.Cells.Item(86, 10).Value = 3687.182  # J86 (was 3627.2727)
$ws.Cells.Item(86, 11).Value = 1110.7333  # K86 (was 1335.125)
$ws.Cells.Item(86, 12).Value = 3687.182  # L86 (was 3627.2727)
$ws.Cells.Item(86, 13).Value = 12.2666999999999  # M86 (was -212.125)
$ws.Cells.Item(86, 14).Value = -5933.182  # N86 (was -5873.2727)
$ws.Cells.Item(87, 8).Value = 50000  # H87 (was 74000)
$ws.Cells.Item(87, 10).Value = 50000  # J87 (was 74000)
$ws.Cells.Item(87, 12).Value = 50000  # L87 (was 74000)
$ws.Cells.Item(87, 14).Value = -52496  # N87 (was -76496)
$ws.Cells.Item(89, 8).Value = 2200.7693  # H89 (was 2662.158)
$ws.Cells.Item(89, 9).Value = 1110.7333  # I89 (was 1335.125)
$ws.Cells.Item(89, 10).Value = 3687.182  # J89 (was 3627.2727)
$ws.Cells.Item(89, 11).Value = 5553.6665  # K89 (was 6675.625)
$ws.Cells.Item(89, 12).Value = 18435.91  # L89 (was 18136.3635)
$ws.Cells.Item(89, 13).Value = 62.33349999999973  # M89 (was -1059.625)
$ws.Cells.Item(89, 14).Value = -29667.91  # N89 (was -29368.3635)
$ws.Cells.Item(90, 8).Value = 50000  # H90 (was 74000)
$ws.Cells.Item(90, 10).Value = 50000  # J90 (was 74000)
$ws.Cells.Item(90, 12).Value = 150000  # L90 (was 222000)
$ws.Cells.Item(90, 14).Value = -162480  # N90 (was -234480)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 27244.266  # H32 (was 25692.807)
$ws.Cells.Item(32, 9).Value = 20594.5  # I32 (was 19804.13)
$ws.Cells.Item(32, 10).Value = 67142.86  # J32 (was 50318.184)
$ws.Cells.Item(32, 11).Value = 20594.5  # K32 (was 19804.13)
$ws.Cells.Item(32, 12).Value = 67142.86  # L32 (was 50318.184)
$ws.Cells.Item(32, 13).Value = -20307.5  # M32 (was -19517.13)
$ws.Cells.Item(32, 14).Value = -67716.86  # N32 (was -50892.184)
$ws.Cells.Item(45, 8).Value = 1138.7059  # H45 (was 1077.5238)
$ws.Cells.Item(45, 9).Value = 1014.6667  # I45 (was 917.55554)
$ws.Cells.Item(45, 10).Value = 1206.3636  # J45 (was 1197.5)
$ws.Cells.Item(45, 11).Value = 1014.6667  # K45 (was 917.55554)
$ws.Cells.Item(45, 12).Value = 1206.3636  # L45 (was 1197.5)
$ws.Cells.Item(45, 13).Value = -637.6667  # M45 (was -540.55554)
$ws.Cells.Item(45, 14).Value = -1960.3636  # N45 (was -1951.5)
$ws.Cells.Item(61, 8).Value = 2711499  # H61 (was 2925548.2)
$ws.Cells.Item(61, 9).Value = 3705037.8  # I61 (was 3969679.8)
$ws.Cells.Item(61, 10).Value = 1847.909  # J61 (was 1979.9)
$ws.Cells.Item(61, 11).Value = 3705037.8  # K61 (was 3969679.8)
$ws.Cells.Item(61, 12).Value = 1847.909  # L61 (was 1979.9)
$ws.Cells.Item(61, 13).Value = -3704825.8  # M61 (was -3969467.8)
$ws.Cells.Item(61, 14).Value = -2271.909  # N61 (was -2403.9)
$ws.Cells.Item(74, 8).Value = 9807949  # H74 (was 9263117)
$ws.Cells.Item(74, 9).Value = 16129817  # I74 (was 16129860)
$ws.Cells.Item(74, 10).Value = 9052.950000000001  # J74 (was 7942.304)
$ws.Cells.Item(74, 11).Value = 16129817  # K74 (was 16129860)
$ws.Cells.Item(74, 12).Value = 9052.950000000001  # L74 (was 7942.304)
$ws.Cells.Item(74, 13).Value = -16128943  # M74 (was -16128986)
$ws.Cells.Item(74, 14).Value = -10800.95  # N74 (was -9690.304)
$ws.Cells.Item(77, 8).Value = 9807949  # H77 (was 9263117)
$ws.Cells.Item(77, 9).Value = 16129817  # I77 (was 16129860)
$ws.Cells.Item(77, 10).Value = 9052.950000000001  # J77 (was 7942.304)
$ws.Cells.Item(77, 11).Value = 80649085  # K77 (was 80649300)
$ws.Cells.Item(77, 12).Value = 45264.75  # L77 (was 39711.52)
$ws.Cells.Item(77, 13).Value = -80644717  # M77 (was -80644932)
$ws.Cells.Item(77, 14).Value = -54000.75  # N77 (was -48447.52)
$ws.Cells.Item(136, 8).Value = 2711499  # H136 (was 2925548.2)
$ws.Cells.Item(136, 9).Value = 3705037.8  # I136 (was 3969679.8)
$ws.Cells.Item(136, 10).Value = 1847.909  # J136 (was 1979.9)
$ws.Cells.Item(136, 11).Value = 11115113.4  # K136 (was 11909039.4)
$ws.Cells.Item(136, 12).Value = 5543.727000000001  # L136 (was 5939.700000000001)
$ws.Cells.Item(136, 13).Value = -11112563.4  # M136 (was -11906489.4)
$ws.Cells.Item(136, 14).Value = -10643.727  # N136 (was -11039.7)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 18244.5  # H31 (was 14709.14)
$ws.Cells.Item(31, 9).Value = 37314.285  # I31 (was 22946.61)
$ws.Cells.Item(31, 10).Value = 7976.154  # J31 (was 7692.037)
$ws.Cells.Item(31, 11).Value = 37314.285  # K31 (was 22946.61)
$ws.Cells.Item(31, 12).Value = 7976.154  # L31 (was 7692.037)
$ws.Cells.Item(31, 13).Value = -37019.285  # M31 (was -22651.61)
$ws.Cells.Item(31, 14).Value = -8566.154  # N31 (was -8282.037)
$ws.Cells.Item(34, 8).Value = 18244.5  # H34 (was 14709.14)
$ws.Cells.Item(34, 9).Value = 37314.285  # I34 (was 22946.61)
$ws.Cells.Item(34, 10).Value = 7976.154  # J34 (was 7692.037)
$ws.Cells.Item(34, 11).Value = 37314.285  # K34 (was 22946.61)
$ws.Cells.Item(34, 12).Value = 7976.154  # L34 (was 7692.037)
$ws.Cells.Item(34, 13).Value = -37112.285  # M34 (was -22744.61)
$ws.Cells.Item(34, 14).Value = -8380.154  # N34 (was -8096.037)
$ws.Cells.Item(58, 8).Value = 2573.204  # H58 (was 3310.2458)
$ws.Cells.Item(58, 9).Value = 1078.52  # I58 (was 1261.7222)
$ws.Cells.Item(58, 10).Value = 4130.1665  # J58 (was 4167.7676)
$ws.Cells.Item(58, 11).Value = 1078.52  # K58 (was 1261.7222)
$ws.Cells.Item(58, 12).Value = 4130.1665  # L58 (was 4167.7676)
$ws.Cells.Item(58, 13).Value = -875.52  # M58 (was -1058.7222)
$ws.Cells.Item(58, 14).Value = -4536.1665  # N58 (was -4573.7676)
$ws.Cells.Item(136, 8).Value = 2573.204  # H136 (was 3310.2458)
$ws.Cells.Item(136, 9).Value = 1078.52  # I136 (was 1261.7222)
$ws.Cells.Item(136, 10).Value = 4130.1665  # J136 (was 4167.7676)
$ws.Cells.Item(136, 11).Value = 3235.56  # K136 (was 3785.1666)
$ws.Cells.Item(136, 12).Value = 12390.4995  # L136 (was 12503.3028)
$ws.Cells.Item(136, 13).Value = -685.5599999999999  # M136 (was -1235.1666)
$ws.Cells.Item(136, 14).Value = -17490.4995  # N136 (was -17603.3028)
$ws.Cells.Item(141, 8).Value = 38043  # H141 (was 41381.832)
$ws.Cells.Item(141, 10).Value = 41358.832  # J141 (was 46028.6)
$ws.Cells.Item(141, 12).Value = 41358.832  # L141 (was 46028.6)
$ws.Cells.Item(141, 14).Value = -51718.832  # N141 (was -56388.6)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 142.88889  # H10 (was 69.14286)
$ws.Cells.Item(10, 9).Value = 123  # I10 (was 69.14286)
$ws.Cells.Item(10, 10).Value = 302  # J10 (was 0)
$ws.Cells.Item(10, 11).Value = 369  # K10 (was 207.42858)
$ws.Cells.Item(10, 12).Value = 906  # L10 (was 0)
$ws.Cells.Item(10, 13).Value = -230  # M10 (was -68.42858000000001)
$ws.Cells.Item(10, 14).Value = -1184  # N10 (was None)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 109402.4  # H80 (was 109414.4)
$ws.Cells.Item(80, 9).Value = 3070  # I80 (was 3095.7144)
$ws.Cells.Item(80, 11).Value = 3070  # K80 (was 3095.7144)
$ws.Cells.Item(80, 13).Value = -2072  # M80 (was -2097.7144)
$ws.Cells.Item(83, 8).Value = 109402.4  # H83 (was 109414.4)
$ws.Cells.Item(83, 9).Value = 3070  # I83 (was 3095.7144)
$ws.Cells.Item(83, 11).Value = 15350  # K83 (was 15478.572)
$ws.Cells.Item(83, 13).Value = -10358  # M83 (was -10486.572)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1649.8334  # H61 (was 1475)
$ws.Cells.Item(61, 9).Value = 1439.8  # I61 (was 1475)
$ws.Cells.Item(61, 10).Value = 2700  # J61 (was 0)
$ws.Cells.Item(61, 11).Value = 1439.8  # K61 (was 1475)
$ws.Cells.Item(61, 12).Value = 2700  # L61 (was 0)
$ws.Cells.Item(61, 13).Value = -1237.8  # M61 (was -1273)
$ws.Cells.Item(61, 14).Value = -3104  # N61 (was None)
$ws.Cells.Item(113, 8).Value = 1649.8334  # H113 (was 1475)
$ws.Cells.Item(113, 9).Value = 1439.8  # I113 (was 1475)
$ws.Cells.Item(113, 10).Value = 2700  # J113 (was 0)
$ws.Cells.Item(113, 11).Value = 1439.8  # K113 (was 1475)
$ws.Cells.Item(113, 12).Value = 2700  # L113 (was 0)
$ws.Cells.Item(113, 13).Value = 730.2  # M113 (was 695)
$ws.Cells.Item(113, 14).Value = -7040  # N113 (was None)
$ws.Cells.Item(122, 8).Value = 2083.5107  # H122 (was 2089.4)
$ws.Cells.Item(122, 9).Value = 2004.7941  # I122 (was 2004.8485)
$ws.Cells.Item(122, 10).Value = 2289.3845  # J122 (was 2321.9167)
$ws.Cells.Item(122, 11).Value = 6014.3823  # K122 (was 6014.5455)
$ws.Cells.Item(122, 12).Value = 6868.1535  # L122 (was 6965.750100000001)
$ws.Cells.Item(122, 13).Value = -3564.3823  # M122 (was -3564.5455)
$ws.Cells.Item(122, 14).Value = -11768.1535  # N122 (was -11865.7501)
